# Auto-generated edit script: update computed market-price columns (H-N)
# across the 8 crafting-leve worksheets, per scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4487.8887
$ws.Range("I64").Value = 3532.3333
$ws.Range("K64").Value = 3532.3333
$ws.Range("M64").Value = -3284.3333

$ws.Range("H67").Value = 4487.8887
$ws.Range("I67").Value = 3532.3333
$ws.Range("K67").Value = 3532.3333
$ws.Range("M67").Value = -2674.3333

$ws.Range("H70").Value = 253299.5
$ws.Range("J70").Value = 336066.34
$ws.Range("L70").Value = 1008199.02
$ws.Range("N70").Value = -1008739.02

$ws.Range("H73").Value = 253299.5
$ws.Range("J73").Value = 336066.34
$ws.Range("L73").Value = 1008199.02
$ws.Range("N73").Value = -1010071.02

$ws.Range("H132").Value = 2438.5625
$ws.Range("I132").Value = 2087.6428
$ws.Range("K132").Value = 6262.928400000001
$ws.Range("M132").Value = -3732.928400000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3671291.5
$ws.Range("I32").Value = 3504856.2
$ws.Range("K32").Value = 3504856.2
$ws.Range("M32").Value = -3504569.2

$ws.Range("H46").Value = 14750
$ws.Range("J46").Value = 5000
$ws.Range("L46").Value = 5000
$ws.Range("N46").Value = -5638

$ws.Range("H63").Value = 14000.167
$ws.Range("I63").Value = 16100.8
$ws.Range("K63").Value = 16100.8
$ws.Range("M63").Value = -15414.8

$ws.Range("H66").Value = 14000.167
$ws.Range("I66").Value = 16100.8
$ws.Range("K66").Value = 80504
$ws.Range("M66").Value = -77072

$ws.Range("H97").Value = 1485
$ws.Range("I97").Value = 1485
$ws.Range("J97").Value = 1485
$ws.Range("K97").Value = 1485
$ws.Range("L97").Value = 1485
$ws.Range("M97").Value = -989
$ws.Range("N97").Value = -2477

$ws.Range("H132").Value = 1507
$ws.Range("I132").Value = 1500
$ws.Range("K132").Value = 4500
$ws.Range("M132").Value = -1970

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1149
$ws.Range("I94").Value = 1149
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1149
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -698
$ws.Range("N94").ClearContents()

$ws.Range("H105").Value = 2077.6
$ws.Range("I105").Value = 2030.6666
$ws.Range("K105").Value = 2030.6666
$ws.Range("M105").Value = -283.6666

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 25.1
$ws.Range("I7").Value = 32.4
$ws.Range("J7").Value = 17.8
$ws.Range("K7").Value = 32.4
$ws.Range("L7").Value = 17.8
$ws.Range("M7").Value = 80.59999999999999
$ws.Range("N7").Value = -243.8

$ws.Range("H31").Value = 2610.2727
$ws.Range("I31").Value = 2471.3
$ws.Range("K31").Value = 2471.3
$ws.Range("M31").Value = -2176.3

$ws.Range("H34").Value = 2610.2727
$ws.Range("I34").Value = 2471.3
$ws.Range("K34").Value = 2471.3
$ws.Range("M34").Value = -2269.3

$ws.Range("H54").Value = 32721.666
$ws.Range("I54").Value = 29083
$ws.Range("K54").Value = 29083
$ws.Range("M54").Value = -28425

$ws.Range("H58").Value = 1642
$ws.Range("I58").Value = 1337.5454
$ws.Range("K58").Value = 1337.5454
$ws.Range("M58").Value = -1134.5454

$ws.Range("H86").Value = 6380.1763
$ws.Range("I86").Value = 6404.8125
$ws.Range("K86").Value = 6404.8125
$ws.Range("M86").Value = -5281.8125

$ws.Range("H89").Value = 6380.1763
$ws.Range("I89").Value = 6404.8125
$ws.Range("K89").Value = 32024.0625
$ws.Range("M89").Value = -26408.0625

$ws.Range("H94").Value = 7110.5
$ws.Range("I94").Value = 7110.5
$ws.Range("K94").Value = 7110.5
$ws.Range("M94").Value = -6659.5

$ws.Range("H99").Value = 1975.7646
$ws.Range("I99").Value = 1912.6
$ws.Range("J99").Value = 2449.5
$ws.Range("K99").Value = 1912.6
$ws.Range("L99").Value = 2449.5
$ws.Range("M99").Value = -414.5999999999999
$ws.Range("N99").Value = -5445.5

$ws.Range("H126").Value = 1975.7646
$ws.Range("I126").Value = 1912.6
$ws.Range("J126").Value = 2449.5
$ws.Range("K126").Value = 5737.799999999999
$ws.Range("L126").Value = 7348.5
$ws.Range("M126").Value = -3267.799999999999
$ws.Range("N126").Value = -12288.5

$ws.Range("H136").Value = 1642
$ws.Range("I136").Value = 1337.5454
$ws.Range("K136").Value = 4012.6362
$ws.Range("M136").Value = -1462.6362

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 127.1
$ws.Range("I2").Value = 7.5
$ws.Range("K2").Value = 45
$ws.Range("M2").Value = 68

$ws.Range("H38").Value = 39.88889
$ws.Range("I38").Value = 41.125
$ws.Range("J38").Value = 30
$ws.Range("K38").Value = 123.375
$ws.Range("L38").Value = 90
$ws.Range("M38").Value = 223.625
$ws.Range("N38").Value = -784

$ws.Range("H41").Value = 1799.5
$ws.Range("I41").Value = 1799.5
$ws.Range("K41").Value = 5398.5
$ws.Range("M41").Value = -5060.5

$ws.Range("H51").Value = 701
$ws.Range("I51").Value = 668
$ws.Range("K51").Value = 2004
$ws.Range("M51").Value = -1544

$ws.Range("H59").Value = 3111
$ws.Range("J59").Value = 3111
$ws.Range("L59").Value = 9333
$ws.Range("N59").Value = -10413

$ws.Range("H81").Value = 1833.5
$ws.Range("I81").Value = 1267
$ws.Range("K81").Value = 3801
$ws.Range("M81").Value = -2678

$ws.Range("H84").Value = 1833.5
$ws.Range("I84").Value = 1267
$ws.Range("K84").Value = 11403
$ws.Range("M84").Value = -5787

$ws.Range("H94").Value = 16337.833
$ws.Range("J94").Value = 16337.833
$ws.Range("L94").Value = 49013.499
$ws.Range("N94").Value = -50365.499

$ws.Range("H132").Value = 4497.8335
$ws.Range("J132").Value = 4745
$ws.Range("L132").Value = 42705
$ws.Range("N132").Value = -47765

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 20004832
$ws.Range("J36").Value = 26671110
$ws.Range("L36").Value = 26671110
$ws.Range("N36").Value = -26672080

$ws.Range("H97").Value = 600
$ws.Range("I97").Value = 200
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 200
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = 296
$ws.Range("N97").Value = -1992

$ws.Range("H107").Value = 2100.5
$ws.Range("I107").Value = 844
$ws.Range("J107").Value = 3859.6
$ws.Range("K107").Value = 844
$ws.Range("L107").Value = 3859.6
$ws.Range("M107").Value = 1076
$ws.Range("N107").Value = -7699.6

$ws.Range("H113").Value = 399
$ws.Range("I113").Value = 399
$ws.Range("K113").Value = 399
$ws.Range("M113").Value = 1771

$ws.Range("H126").Value = 2399.6
$ws.Range("I126").Value = 1999.5
$ws.Range("K126").Value = 5998.5
$ws.Range("M126").Value = -3528.5

$ws.Range("H132").Value = 19800
$ws.Range("I132").Value = 19800
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 59400
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -56870
$ws.Range("N132").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1774.875
$ws.Range("J46").Value = 1975
$ws.Range("L46").Value = 1975
$ws.Range("N46").Value = -2351

$ws.Range("H82").Value = 700.75
$ws.Range("I82").Value = 459.42856
$ws.Range("J82").Value = 1038.6
$ws.Range("K82").Value = 459.42856
$ws.Range("L82").Value = 1038.6
$ws.Range("M82").Value = -98.42856
$ws.Range("N82").Value = -1760.6

$ws.Range("H85").Value = 700.75
$ws.Range("I85").Value = 459.42856
$ws.Range("J85").Value = 1038.6
$ws.Range("K85").Value = 459.42856
$ws.Range("L85").Value = 1038.6
$ws.Range("M85").Value = 788.5714399999999
$ws.Range("N85").Value = -3534.6

$ws.Range("H93").Value = 1243.8334
$ws.Range("I93").Value = 1054.5
$ws.Range("K93").Value = 1054.5
$ws.Range("M93").Value = 193.5

$ws.Range("H100").Value = 2998.2354
$ws.Range("I100").Value = 2871.4
$ws.Range("K100").Value = 2871.4
$ws.Range("M100").Value = -2330.4

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()

$ws.Range("H29").Value = 6999
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

$ws.Range("H41").Value = 15398.429
$ws.Range("I41").Value = 11325
$ws.Range("K41").Value = 11325
$ws.Range("M41").Value = -10935

$ws.Range("H75").Value = 29559
$ws.Range("I75").Value = 29559
$ws.Range("K75").Value = 29559
$ws.Range("M75").Value = -28623

$ws.Range("H78").Value = 29559
$ws.Range("I78").Value = 29559
$ws.Range("K78").Value = 88677
$ws.Range("M78").Value = -83997

$ws.Range("H107").Value = 2252.3914
$ws.Range("I107").Value = 2238.7222
$ws.Range("J107").Value = 2301.6
$ws.Range("K107").Value = 6716.1666
$ws.Range("L107").Value = 6904.799999999999
$ws.Range("M107").Value = -4796.1666
$ws.Range("N107").Value = -10744.8

$ws.Range("H126").Value = 3792.2
$ws.Range("I126").Value = 2531
$ws.Range("K126").Value = 7593
$ws.Range("M126").Value = -5123
